$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-point the "CauseOfDeathURI*" cells at the new ICD-11 "release" URL scheme and
# turn them into real hyperlinks (this is what Excel does automatically: it adds a
# "Hyperlink" cell style/font and records the link target in the worksheet rels).

$links = @(
    @{ Cell = "I2";  Url = "http://id.who.int/icd/release/11/mms/1334938734/unspecified" },
    @{ Cell = "I3";  Url = "http://id.who.int/icd/release/11/mms/1458683894/unspecified" },
    @{ Cell = "I4";  Url = "http://id.who.int/icd/release/11/mms/1974956233" },
    @{ Cell = "I5";  Url = "http://id.who.int/icd/release/11/mms/1974956233" },

    @{ Cell = "M2";  Url = "http://id.who.int/icd/release/11/mms/316539081/unspecified" },
    @{ Cell = "M3";  Url = "http://id.who.int/icd/release/11/mms/512128824" },
    @{ Cell = "M4";  Url = "http://id.who.int/icd/release/11/mms/2004408087/unspecified" },
    @{ Cell = "M5";  Url = "http://id.who.int/icd/release/11/mms/2004408087/unspecified" },

    @{ Cell = "Q2";  Url = "http://id.who.int/icd/release/11/mms/1917052637/unspecified" },
    @{ Cell = "Q3";  Url = "http://id.who.int/icd/release/11/mms/1997348476/unspecified" },
    @{ Cell = "Q4";  Url = "http://id.who.int/icd/release/11/mms/1178642763" },
    @{ Cell = "Q5";  Url = "http://id.who.int/icd/release/11/mms/1178642763" },

    @{ Cell = "U2";  Url = "http://id.who.int/icd/release/11/mms/761947693/unspecified" },
    @{ Cell = "U3";  Url = "http://id.who.int/icd/release/11/mms/770085732" },

    @{ Cell = "AC3"; Url = "http://id.who.int/icd/release/11/mms/1580466198/unspecified" },
    @{ Cell = "AC4"; Url = "http://id.who.int/icd/release/11/mms/1726201225/unspecified" },
    @{ Cell = "AC5"; Url = "http://id.who.int/icd/release/11/mms/1726201225/unspecified" }
)

foreach ($link in $links) {
    $rng = $ws.Range($link.Cell)
    $ws.Hyperlinks.Add($rng, $link.Url, "", "", $link.Url) | Out-Null
}

# Widen the columns that now hold full URIs / long text so the links are readable,
# matching the column sizing that Excel applied when the sheet was reformatted.
$ws.Range("H1").EntireColumn.ColumnWidth = 15.28515625
$ws.Range("I1").EntireColumn.ColumnWidth = 58.85546875
$ws.Range("J1").EntireColumn.ColumnWidth = 22.28515625
$ws.Range("M1").EntireColumn.ColumnWidth = 67.42578125
$ws.Range("Q1").EntireColumn.ColumnWidth = 69.140625
$ws.Range("U1").EntireColumn.ColumnWidth = 57.7109375
$ws.Range("AB1").EntireColumn.ColumnWidth = 23.28515625
$ws.Range("AC1").EntireColumn.ColumnWidth = 71.85546875

# Scroll the view over towards the newly widened columns and select where editing left off.
$ws.Application.ActiveWindow.ScrollColumn = 22
$ws.Range("AC6").Select()
